$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated DAMSLTag (column I) / DialogAct (column J) values following
# re-run of SGNN dialog act annotation.
$updates = @(
    @{ Row = 5;   I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 12;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 14;  I = "%";  J = "Uninterpretable" },
    @{ Row = 16;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 20;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 25;  I = "ba"; J = "Appreciation" },
    @{ Row = 27;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 34;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 35;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 38;  I = "aa"; J = "Agree/Accept" },
    @{ Row = 39;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 45;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 50;  I = "aa"; J = "Agree/Accept" },
    @{ Row = 51;  I = "b";  J = "Acknowledge (Backchannel)" },
    @{ Row = 69;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 74;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 75;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 80;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 81;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 93;  I = "b";  J = "Acknowledge (Backchannel)" },
    @{ Row = 98;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 103; I = "sv"; J = "Statement-opinion" },
    @{ Row = 111; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 120; I = "sv"; J = "Statement-opinion" },
    @{ Row = 121; I = "aa"; J = "Agree/Accept" },
    @{ Row = 132; I = "sd"; J = "Statement-non-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
